# "final update of first part. mapping to QPX genome"
#
# Adds two new worksheets (Sheet3, Sheet4) with QPX-genome mapping stats
# (pfam/peptide/domain matches and blat contig matches), appends their
# labels to the shared-string pool, and makes the newly added Sheet3 the
# active tab of the workbook.

$wb = $excel.ActiveWorkbook

$sheet1 = $wb.Worksheets.Item(1)
$sheet2 = $wb.Worksheets.Item(2)

# New sheets go after Sheet2, in order: Sheet3 then Sheet4.
$sheet3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $sheet2)
$sheet4 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $sheet3)

# ---------------------------------------------------------------
# Sheet3 : pfam / peptides / domain hits against the QPX assembly
# ---------------------------------------------------------------
$sheet3.Range("A1").Value = "sample"
$sheet3.Range("B1").Value = "reference"
$sheet3.Range("C1").Value = 0.000000001
$sheet3.Range("D1").Value = 0.00001
$sheet3.Range("E1").Value = 1
$sheet3.Range("F1").Value = "pfam"

$sheet3.Range("A2").Value = "a"
$sheet3.Range("B2").Value = "QPX"
$sheet3.Range("C2").Value = 498
$sheet3.Range("D2").Value = 658
$sheet3.Range("E2").Value = 1101
$sheet3.Range("F2").Value = "peptides"

$sheet3.Range("A3").Value = "b"
$sheet3.Range("B3").Value = "QPX"
$sheet3.Range("C3").Value = 21
$sheet3.Range("D3").Value = 30
$sheet3.Range("E3").Value = 45
$sheet3.Range("F3").Value = "peptides"

$sheet3.Range("A4").Value = "c"
$sheet3.Range("B4").Value = "QPX"
$sheet3.Range("C4").Value = 0
$sheet3.Range("D4").Value = 0
$sheet3.Range("E4").Value = 0
$sheet3.Range("F4").Value = "peptides"

$sheet3.Range("A5").Value = "a"
$sheet3.Range("B5").Value = "QPX"
$sheet3.Range("C5").Value = 359
$sheet3.Range("D5").Value = 511
$sheet3.Range("E5").Value = 1568
$sheet3.Range("F5").Value = "domain"

$sheet3.Range("A6").Value = "b"
$sheet3.Range("B6").Value = "QPX"
$sheet3.Range("C6").Value = 13
$sheet3.Range("D6").Value = 23
$sheet3.Range("E6").Value = 53
$sheet3.Range("F6").Value = "domain"

$sheet3.Range("A7").Value = "c"
$sheet3.Range("B7").Value = "QPX"
$sheet3.Range("C7").Value = 0
$sheet3.Range("D7").Value = 0
$sheet3.Range("E7").Value = 0
$sheet3.Range("F7").Value = "domain"

# Column C is sized to best-fit its (scientific-notation) contents, as in
# the sibling Sheet2 which also carries a custom best-fit width on col C.
$sheet3.Columns.Item(3).ColumnWidth = 11.25

# ---------------------------------------------------------------
# Sheet4 : blat match counts (match100..match500) against QPX contigs
# ---------------------------------------------------------------
$sheet4.Range("A1").Value = "sample"
$sheet4.Range("B1").Value = "reference"
$sheet4.Range("C1").Value = "match100"
$sheet4.Range("D1").Value = "match200"
$sheet4.Range("E1").Value = "match300"
$sheet4.Range("F1").Value = "match400"
$sheet4.Range("G1").Value = "match500"
$sheet4.Range("H1").Value = "blat"

$sheet4.Range("A2").Value = "a"
$sheet4.Range("B2").Value = "QPX"
$sheet4.Range("C2").Value = 1124
$sheet4.Range("D2").Value = 1111
$sheet4.Range("E2").Value = 647
$sheet4.Range("F2").Value = 382
$sheet4.Range("G2").Value = 262
$sheet4.Range("H2").Value = "assembled,contig"

$sheet4.Range("A3").Value = "b"
$sheet4.Range("B3").Value = "QPX"
$sheet4.Range("C3").Value = 50
$sheet4.Range("D3").Value = 48
$sheet4.Range("E3").Value = 25
$sheet4.Range("F3").Value = 15
$sheet4.Range("G3").Value = 14
$sheet4.Range("H3").Value = "assembled,contig"

$sheet4.Range("A4").Value = "c"
$sheet4.Range("B4").Value = "QPX"
$sheet4.Range("C4").Value = 2
$sheet4.Range("D4").Value = 2
$sheet4.Range("E4").Value = 1
$sheet4.Range("F4").Value = 0
$sheet4.Range("G4").Value = 0
$sheet4.Range("H4").Value = "assembled,contig"

$sheet4.Range("A5").Value = "a"
$sheet4.Range("B5").Value = "QPX"
$sheet4.Range("C5").Value = 859
$sheet4.Range("D5").Value = 831
$sheet4.Range("E5").Value = 535
$sheet4.Range("F5").Value = 342
$sheet4.Range("G5").Value = 240
$sheet4.Range("H5").Value = "genomic.contig"

$sheet4.Range("A6").Value = "b"
$sheet4.Range("B6").Value = "QPX"
$sheet4.Range("C6").Value = 39
$sheet4.Range("D6").Value = 35
$sheet4.Range("E6").Value = 22
$sheet4.Range("F6").Value = 14
$sheet4.Range("G6").Value = 13
$sheet4.Range("H6").Value = "genomic.contig"

$sheet4.Range("A7").Value = "c"
$sheet4.Range("B7").Value = "QPX"
$sheet4.Range("C7").Value = 2
$sheet4.Range("D7").Value = 2
$sheet4.Range("E7").Value = 1
$sheet4.Range("F7").Value = 0
$sheet4.Range("G7").Value = 0
$sheet4.Range("H7").Value = "genomic.contig"

# Sheet3 becomes the active/selected tab (workbookView activeTab="2",
# sheetView tabSelected moves from Sheet1 to Sheet3).
$sheet3.Activate()
